$p = $ppt.ActivePresentation

# --- 1. On the closing "Thank You" slide (still slide 10 at this point),
#        shrink the title text to 35pt. The title textbox auto-fits its
#        text (wrap="none" + spAutoFit), so shrink/recenter the shape to
#        match what PowerPoint would compute for the smaller font. ---
$thankYouSlide = $p.Slides.Item($p.Slides.Count)
$titleShape = $thankYouSlide.Shapes.Item("TextBox 1")
$titleShape.TextFrame.TextRange.Font.Size = 35
$titleShape.Left = 362.27752685546875
$titleShape.Top = 144.00003051757812
$titleShape.Width = 211.4449920654297
$titleShape.Height = 49.68048858642578

# --- 2. Insert a new "Title and Content" slide at position 3 (right after
#        the Agenda slide, before the Simplicity slide). ---
$titleAndContentLayout = $p.SlideMaster.CustomLayouts.Item(2)
$newSlide = $p.Slides.AddSlide(3, $titleAndContentLayout)
$newSlide.Shapes.Item(1).TextFrame.TextRange.Text = "Hello"
$newSlide.Shapes.Item(2).TextFrame.TextRange.Text = "Bob"
